$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.713.12'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.836.96'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.37'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4683'
$ws.Range('E7').Value = '  +3.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3615'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07160'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9329'
$ws.Range('E10').Value = '  +4.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.52'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07665'
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.868.16'
$ws.Range('E13').Value = '  +3.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.270'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.371'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.92'
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008561'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.738.17'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.29'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.59'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.917'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.90'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.98'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.010'
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.88'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.910'
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08831'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('B31').Value = 'HuobiToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.154'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.836'
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.184'
$ws.Range('E33').Value = '  +6.64%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7412'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.449'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.083'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.964'
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01926'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05153'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.924'
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5077'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1502'
$ws.Range('E42').Value = '  -0.67%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.127'
$ws.Range('E43').Value = '  +1.36%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4674'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.007'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.19'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.45'
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.580'
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06034'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.33'
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.07'
$ws.Range('E51').Value = '  -0.16%  '
